$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like numbers but must remain stored as literal text
# (matching the source data's inlineStr representation, e.g. '1.00' not 1).
# Force text format, assign, then restore the default 'Normal' style so no
# extra formatting is visibly applied to the cell.
$textForcedCells = @(
    @{Addr = 'D5'; Value = '585.41'},
    @{Addr = 'D6'; Value = '148.07'},
    @{Addr = 'D10'; Value = '5.73'},
    @{Addr = 'D13'; Value = '27.34'},
    @{Addr = 'D18'; Value = '12.11'},
    @{Addr = 'D19'; Value = '4.65'},
    @{Addr = 'D20'; Value = '343.93'},
    @{Addr = 'D21'; Value = '6.80'},
    @{Addr = 'D23'; Value = '66.82'},
    @{Addr = 'D24'; Value = '1.69'},
    @{Addr = 'D25'; Value = '9.05'},
    @{Addr = 'D27'; Value = '552.64'},
    @{Addr = 'D28'; Value = '7.98'},
    @{Addr = 'D29'; Value = '0.161'},
    @{Addr = 'D35'; Value = '164.99'},
    @{Addr = 'D36'; Value = '0.411'},
    @{Addr = 'D37'; Value = '1.00'},
    @{Addr = 'D40'; Value = '0.999'},
    @{Addr = 'D41'; Value = '165.51'},
    @{Addr = 'D42'; Value = '39.63'},
    @{Addr = 'D43'; Value = '3.93'},
    @{Addr = 'D44'; Value = '0.0585'},
    @{Addr = 'D45'; Value = '22.50'},
    @{Addr = 'D46'; Value = '0.628'},
    @{Addr = 'D48'; Value = '0.0246'},
    @{Addr = 'D50'; Value = '18.89'}
)
foreach ($item in $textForcedCells) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

# Remaining plain text / percentage / url / label cell updates
$plainCells = @(
    @{Addr = 'D2'; Value = '63.150.01'},
    @{Addr = 'E2'; Value = '  -2.04%  '},
    @{Addr = 'D3'; Value = '2.574.59'},
    @{Addr = 'E3'; Value = '  -3.21%  '},
    @{Addr = 'E4'; Value = '  +0.06%  '},
    @{Addr = 'E5'; Value = '  -3.84%  '},
    @{Addr = 'E6'; Value = '  -3.43%  '},
    @{Addr = 'E7'; Value = '  +0.06%  '},
    @{Addr = 'E8'; Value = '  -1.30%  '},
    @{Addr = 'E9'; Value = '  -1.07%  '},
    @{Addr = 'E10'; Value = '  +2.14%  '},
    @{Addr = 'E11'; Value = '  -1.73%  '},
    @{Addr = 'E12'; Value = '  -0.89%  '},
    @{Addr = 'E13'; Value = '  -3.08%  '},
    @{Addr = 'D14'; Value = '3.037.96'},
    @{Addr = 'E14'; Value = '  -3.15%  '},
    @{Addr = 'D15'; Value = '63.095.76'},
    @{Addr = 'E15'; Value = '  -1.87%  '},
    @{Addr = 'E16'; Value = '  +2.48%  '},
    @{Addr = 'D17'; Value = '2.562.47'},
    @{Addr = 'E17'; Value = '  -3.51%  '},
    @{Addr = 'E18'; Value = '  -0.47%  '},
    @{Addr = 'E20'; Value = '  -2.00%  '},
    @{Addr = 'E21'; Value = '  -2.10%  '},
    @{Addr = 'E23'; Value = '  +0.01%  '},
    @{Addr = 'E24'; Value = '  -3.24%  '},
    @{Addr = 'E25'; Value = '  -3.84%  '},
    @{Addr = 'E26'; Value = '  -4.26%  '},
    @{Addr = 'E27'; Value = '  -0.53%  '},
    @{Addr = 'B28'; Value = 'Aptos'},
    @{Addr = 'C28'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'},
    @{Addr = 'E28'; Value = '  -2.78%  '},
    @{Addr = 'B29'; Value = 'Kaspa'},
    @{Addr = 'C29'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'},
    @{Addr = 'E29'; Value = '  -1.59%  '},
    @{Addr = 'E30'; Value = '  +0.11%  '},
    @{Addr = 'E31'; Value = '  -2.59%  '},
    @{Addr = 'D32'; Value = '0.0₃0846'},
    @{Addr = 'E32'; Value = '  -2.41%  '},
    @{Addr = 'E33'; Value = '  -2.14%  '},
    @{Addr = 'E34'; Value = '  -3.96%  '},
    @{Addr = 'E35'; Value = '  -2.12%  '},
    @{Addr = 'E36'; Value = '  +0.42%  '},
    @{Addr = 'E37'; Value = '  -0.02%  '},
    @{Addr = 'E39'; Value = '  -4.95%  '},
    @{Addr = 'E40'; Value = '  +0.01%  '},
    @{Addr = 'E41'; Value = '  -0.88%  '},
    @{Addr = 'E42'; Value = '  -1.54%  '},
    @{Addr = 'E43'; Value = '  +1.86%  '},
    @{Addr = 'E44'; Value = '  +1.23%  '},
    @{Addr = 'E45'; Value = '  +1.78%  '},
    @{Addr = 'E46'; Value = '  -0.75%  '},
    @{Addr = 'E47'; Value = '  +0.61%  '},
    @{Addr = 'E48'; Value = '  -0.31%  '},
    @{Addr = 'E49'; Value = '  -1.01%  '},
    @{Addr = 'E50'; Value = '  -1.18%  '},
    @{Addr = 'D51'; Value = '0.0₆0224'},
    @{Addr = 'E51'; Value = '  +11.24%  '}
)
foreach ($item in $plainCells) {
    $ws.Range($item.Addr).Value = $item.Value
}

